$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Compliance Overview": bump violation counts for RDVI row (9)
# and the Totals row (10) from 2 to 9.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Compliance Overview")
$wsOverview.Range("C9").Value = 9
$wsOverview.Range("D9").Value = 9
$wsOverview.Range("C10").Value = 9
$wsOverview.Range("D10").Value = 9

# ---------------------------------------------------------------------
# Sheet "Compliance Details": expand the RDVI compliance-check detail
# table from 2 data rows to 10 data rows (rows 2-11), inserting the
# additional checks alphabetically while keeping the FAIL styling
# (style copied from the existing row 2) on every row except the new
# "prospectus_80pct_policy" PASS row, which keeps default formatting.
# ---------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("Compliance Details")

# First, move the current row 3 ("max_15pct_illiquid_sai") data down to
# row 6 so the new rows can be inserted above it in alphabetical order.
$wsDetails.Range("A3:G3").Copy()
$wsDetails.Range("A6:G6").PasteSpecial(-4104)
$wsDetails.Range("A3:G3").ClearContents()

# Rename the check on row 2.
$wsDetails.Range("B2").Value = "commodities_check"

# Give rows 3-11 the same red FAIL fill used on row 2 by default; the
# PASS row (7) will have its formatting cleared afterwards.
$wsDetails.Range("A2:G2").Copy()
$wsDetails.Range("A3:G11").PasteSpecial(-4122)

# Row 3: diversification_40act_check
$wsDetails.Range("A3").Value = "RDVI"
$wsDetails.Range("B3").Value = "diversification_40act_check"
$wsDetails.Range("C3").Value = "FAIL"
$wsDetails.Range("D3").Value = "FAIL"
$wsDetails.Range("E3").Value = 1
$wsDetails.Range("F3").Value = 1
$wsDetails.Range("G3").Value = "NO"

# Row 4: diversification_IRS_check
$wsDetails.Range("A4").Value = "RDVI"
$wsDetails.Range("B4").Value = "diversification_IRS_check"
$wsDetails.Range("C4").Value = "FAIL"
$wsDetails.Range("D4").Value = "FAIL"
$wsDetails.Range("E4").Value = 1
$wsDetails.Range("F4").Value = 1
$wsDetails.Range("G4").Value = "NO"

# Row 5: gics_compliance
$wsDetails.Range("A5").Value = "RDVI"
$wsDetails.Range("B5").Value = "gics_compliance"
$wsDetails.Range("C5").Value = "FAIL"
$wsDetails.Range("D5").Value = "FAIL"
$wsDetails.Range("E5").Value = 1
$wsDetails.Range("F5").Value = 1
$wsDetails.Range("G5").Value = "NO"

# Row 6: max_15pct_illiquid_sai (values re-asserted; formatting/values
# already copied above from the original row 3)
$wsDetails.Range("A6").Value = "RDVI"
$wsDetails.Range("B6").Value = "max_15pct_illiquid_sai"
$wsDetails.Range("C6").Value = "FAIL"
$wsDetails.Range("D6").Value = "FAIL"
$wsDetails.Range("E6").Value = 1
$wsDetails.Range("F6").Value = 1
$wsDetails.Range("G6").Value = "NO"

# Row 7: prospectus_80pct_policy (PASS - no special fill)
$wsDetails.Range("A7").Value = "RDVI"
$wsDetails.Range("B7").Value = "prospectus_80pct_policy"
$wsDetails.Range("C7").Value = "PASS"
$wsDetails.Range("D7").Value = "PASS"
$wsDetails.Range("E7").Value = 0
$wsDetails.Range("F7").Value = 0
$wsDetails.Range("G7").Value = "NO"
$wsDetails.Range("A7:G7").ClearFormats()

# Row 8: real_estate_check
$wsDetails.Range("A8").Value = "RDVI"
$wsDetails.Range("B8").Value = "real_estate_check"
$wsDetails.Range("C8").Value = "FAIL"
$wsDetails.Range("D8").Value = "FAIL"
$wsDetails.Range("E8").Value = 1
$wsDetails.Range("F8").Value = 1
$wsDetails.Range("G8").Value = "NO"

# Row 9: twelve_d1a_other_inv_cos
$wsDetails.Range("A9").Value = "RDVI"
$wsDetails.Range("B9").Value = "twelve_d1a_other_inv_cos"
$wsDetails.Range("C9").Value = "FAIL"
$wsDetails.Range("D9").Value = "FAIL"
$wsDetails.Range("E9").Value = 1
$wsDetails.Range("F9").Value = 1
$wsDetails.Range("G9").Value = "NO"

# Row 10: twelve_d2_insurance_cos
$wsDetails.Range("A10").Value = "RDVI"
$wsDetails.Range("B10").Value = "twelve_d2_insurance_cos"
$wsDetails.Range("C10").Value = "FAIL"
$wsDetails.Range("D10").Value = "FAIL"
$wsDetails.Range("E10").Value = 1
$wsDetails.Range("F10").Value = 1
$wsDetails.Range("G10").Value = "NO"

# Row 11: twelve_d3_sec_biz
$wsDetails.Range("A11").Value = "RDVI"
$wsDetails.Range("B11").Value = "twelve_d3_sec_biz"
$wsDetails.Range("C11").Value = "FAIL"
$wsDetails.Range("D11").Value = "FAIL"
$wsDetails.Range("E11").Value = 1
$wsDetails.Range("F11").Value = 1
$wsDetails.Range("G11").Value = "NO"
